{"js": "// Replace the arithmetic expressions in the 20x5 practice-sheet table with\n// the new set of expressions (text-only change; fonts/sizes/layout untouched).\n// The new values below are taken row-by-row, left-to-right, straight from the\n// target OOXML diff.\nconst newValues = [\n  [\"92+2=\", \"23+70=\", \"0+55=\", \"94-37=\", \"61-53=\"],\n  [\"56+35=\", \"37+40=\", \"19+34=\", \"49-8=\", \"74-8=\"],\n  [\"14+75=\", \"58+29=\", \"59-52=\", \"47-46=\", \"34-32=\"],\n  [\"28+28=\", \"65-58=\", \"4+19=\", \"58+11=\", \"33+46=\"],\n  [\"37-10=\", \"50-33=\", \"11+29=\", \"30+62=\", \"36-17=\"],\n  [\"2+6=\", \"86-42=\", \"9+84=\", \"71+22=\", \"8+86=\"],\n  [\"78+8=\", \"45-2=\", \"61+0=\", \"5-2=\", \"5+78=\"],\n  [\"20+10=\", \"53+28=\", \"75-6=\", \"7+52=\", \"7+78=\"],\n  [\"94-28=\", \"4+46=\", \"85-18=\", \"96-9=\", \"45+47=\"],\n  [\"45+12=\", \"96-47=\", \"83-61=\", \"52+20=\", \"42-15=\"],\n  [\"51+14=\", \"26-11=\", \"29+20=\", \"42+51=\", \"33-23=\"],\n  [\"10+3=\", \"62-28=\", \"42-29=\", \"74-7=\", \"1+25=\"],\n  [\"38-30=\", \"47-47=\", \"55+21=\", \"11+76=\", \"14+62=\"],\n  [\"38+31=\", \"99-42=\", \"87-37=\", \"88-67=\", \"85-50=\"],\n  [\"56-52=\", \"59+6=\", \"40-27=\", \"94-38=\", \"33-28=\"],\n  [\"82-8=\", \"21+73=\", \"9+74=\", \"27-17=\", \"53-0=\"],\n  [\"20+55=\", \"26+53=\", \"48+38=\", \"67-58=\", \"17-9=\"],\n  [\"79-28=\", \"46+49=\", \"69-10=\", \"82-67=\", \"77+8=\"],\n  [\"78-51=\", \"72+18=\", \"33+0=\", \"81-49=\", \"44+18=\"],\n  [\"87-53=\", \"95-63=\", \"62+34=\", \"79-44=\", \"29+56=\"],\n];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nfor (let r = 0; r < newValues.length; r++) {\n  const row = newValues[r];\n  for (let c = 0; c < row.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = row[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the arithmetic expressions in the 20x5 practice-sheet table with\n# the new set of expressions (text-only change; fonts/sizes/layout untouched).\n# The new values below are taken row-by-row, left-to-right, straight from the\n# target OOXML diff.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"92+2=\", \"23+70=\", \"0+55=\", \"94-37=\", \"61-53=\"),\n    @(\"56+35=\", \"37+40=\", \"19+34=\", \"49-8=\", \"74-8=\"),\n    @(\"14+75=\", \"58+29=\", \"59-52=\", \"47-46=\", \"34-32=\"),\n    @(\"28+28=\", \"65-58=\", \"4+19=\", \"58+11=\", \"33+46=\"),\n    @(\"37-10=\", \"50-33=\", \"11+29=\", \"30+62=\", \"36-17=\"),\n    @(\"2+6=\", \"86-42=\", \"9+84=\", \"71+22=\", \"8+86=\"),\n    @(\"78+8=\", \"45-2=\", \"61+0=\", \"5-2=\", \"5+78=\"),\n    @(\"20+10=\", \"53+28=\", \"75-6=\", \"7+52=\", \"7+78=\"),\n    @(\"94-28=\", \"4+46=\", \"85-18=\", \"96-9=\", \"45+47=\"),\n    @(\"45+12=\", \"96-47=\", \"83-61=\", \"52+20=\", \"42-15=\"),\n    @(\"51+14=\", \"26-11=\", \"29+20=\", \"42+51=\", \"33-23=\"),\n    @(\"10+3=\", \"62-28=\", \"42-29=\", \"74-7=\", \"1+25=\"),\n    @(\"38-30=\", \"47-47=\", \"55+21=\", \"11+76=\", \"14+62=\"),\n    @(\"38+31=\", \"99-42=\", \"87-37=\", \"88-67=\", \"85-50=\"),\n    @(\"56-52=\", \"59+6=\", \"40-27=\", \"94-38=\", \"33-28=\"),\n    @(\"82-8=\", \"21+73=\", \"9+74=\", \"27-17=\", \"53-0=\"),\n    @(\"20+55=\", \"26+53=\", \"48+38=\", \"67-58=\", \"17-9=\"),\n    @(\"79-28=\", \"46+49=\", \"69-10=\", \"82-67=\", \"77+8=\"),\n    @(\"78-51=\", \"72+18=\", \"33+0=\", \"81-49=\", \"44+18=\"),\n    @(\"87-53=\", \"95-63=\", \"62+34=\", \"79-44=\", \"29+56=\")\n)\n\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n    $row = $newValues[$r]\n    for ($c = 0; $c -lt $row.Length; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n    }\n}\n"}
